$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "29.580.73"; E = "  -2.26%  " },
    @{ Row = 3; D = "2.003.05"; E = "  -4.03%  " },
    @{ Row = 4; D = "1.012"; E = "  +0.86%  " },
    @{ Row = 5; D = "329.94"; E = "  -3.67%  " },
    @{ Row = 6; D = "1.011"; E = "  +0.83%  " },
    @{ Row = 7; D = "0.5015"; E = "  -4.02%  " },
    @{ Row = 8; D = "0.4222"; E = "  -4.10%  " },
    @{ Row = 9; D = "54.03"; E = "  -0.89%  " },
    @{ Row = 10; D = "0.09007"; E = "  -3.64%  " },
    @{ Row = 11; D = "1.119"; E = "  -4.04%  " },
    @{ Row = 12; D = "23.30"; E = "  -5.89%  " },
    @{ Row = 13; D = "2.003.68"; E = "  -0.01%  " },
    @{ Row = 14; D = "8.049"; E = "  -5.95%  " },
    @{ Row = 15; D = "6.491"; E = "  -5.73%  " },
    @{ Row = 16; D = "1.014"; E = "  +0.91%  " },
    @{ Row = 17; D = "94.42"; E = "  -6.72%  " },
    @{ Row = 18; D = "0.00001113"; E = "  -3.69%  " },
    @{ Row = 19; D = "0.06688"; E = "  +0.36%  " },
    @{ Row = 20; D = "19.69"; E = "  -6.49%  " },
    @{ Row = 21; D = "1.011"; E = "  +0.87%  " },
    @{ Row = 22; D = "5.965"; E = "  -5.52%  " },
    @{ Row = 23; D = "29.613.70"; E = "  -2.13%  " },
    @{ Row = 24; D = "12.01"; E = "  -3.86%  " },
    @{ Row = 25; E = "  -0.10%  " },
    @{ Row = 26; D = "158.92"; E = "  -2.00%  " },
    @{ Row = 27; D = "20.72"; E = "  -4.76%  " },
    @{ Row = 28; D = "6.402"; E = "  -4.02%  " },
    @{ Row = 29; D = "2.299"; E = "  -8.20%  " },
    @{ Row = 30; D = "128.23"; E = "  -3.56%  " },
    @{ Row = 31; D = "1.057"; E = "  -6.41%  " },
    @{ Row = 32; D = "0.09963"; E = "  -4.51%  " },
    @{ Row = 33; D = "1.569"; E = "  -5.58%  " },
    @{ Row = 34; D = "5.842"; E = "  -5.96%  " },
    @{ Row = 35; D = "3.795"; E = "  -1.69%  " },
    @{ Row = 36; D = "0.02475"; E = "  -5.52%  " },
    @{ Row = 37; D = "9.304"; E = "  -8.74%  " },
    @{ Row = 38; D = "1.308"; E = "  -2.70%  " },
    @{ Row = 39; D = "0.06392"; E = "  -6.23%  " },
    @{ Row = 40; D = "0.6561"; E = "  -5.69%  " },
    @{ Row = 41; D = "11.71"; E = "  -6.24%  " },
    @{ Row = 42; D = "0.2057"; E = "  -6.66%  " },
    @{ Row = 43; D = "1.011"; E = "  +0.85%  " },
    @{ Row = 44; D = "0.6359"; E = "  -6.49%  " },
    @{ Row = 45; D = "13.42"; E = "  -5.84%  " },
    @{ Row = 46; D = "2.195"; E = "  -5.40%  " },
    @{ Row = 47; D = "1.304"; E = "  -4.70%  " },
    @{ Row = 48; E = "  -3.22%  " },
    @{ Row = 49; E = "  -2.46%  " },
    @{ Row = 50; D = "0.06991"; E = "  -3.24%  " },
    @{ Row = 51; D = "1.128"; E = "  -6.69%  " }
)

foreach ($item in $updates) {
    if ($item.ContainsKey("D")) {
        $cellD = $ws.Cells.Item($item.Row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $item.D
        $cellD.Style = "Normal"
    }
    if ($item.ContainsKey("E")) {
        $ws.Cells.Item($item.Row, 5).Value = $item.E
    }
}
